$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$ws.Range("D2").Value = 1318
$ws.Range("E2").Value = 74
$ws.Range("F2").Value = 74
$ws.Range("G2").Value = 82
$ws.Range("H2").Value = 62
$ws.Range("I2").Value = 62
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1325
$ws.Range("L2").Value = 690
$ws.Range("M2").Value = 635
$ws.Range("N2").Value = 635
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 28
$ws.Range("Q2").Value = 39
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = -13
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = 38
$ws.Range("V2").Value = 195
$ws.Range("W2").Value = 5.61
$ws.Range("X2").Value = 4.72
$ws.Range("Y2").Value = 10.07
$ws.Range("Z2").Value = 4.76
$ws.Range("AA2").Value = 108.75
$ws.Range("AB2").Value = 2379.16
$ws.Range("AC2").Value = 557
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").Value = 5688
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 100
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").Value = 17.95
$ws.Range("AJ2").Value = 11159460

# Row 3
$ws.Range("D3").Value = 1746
$ws.Range("E3").Value = 122
$ws.Range("F3").Value = 122
$ws.Range("G3").Value = 177
$ws.Range("H3").Value = 141
$ws.Range("I3").Value = 141
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2183
$ws.Range("L3").Value = 588
$ws.Range("M3").Value = 1595
$ws.Range("N3").Value = 1595
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 41
$ws.Range("Q3").Value = 78
$ws.Range("R3").Value = -650
$ws.Range("S3").Value = 646
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 77
$ws.Range("V3").Value = 30
$ws.Range("W3").Value = 6.98
$ws.Range("X3").Value = 8.08
$ws.Range("Y3").Value = 12.65
$ws.Range("Z3").Value = 8.04
$ws.Range("AA3").Value = 36.88
$ws.Range("AB3").Value = 3924.28
$ws.Range("AC3").Value = 1110
$ws.Range("AD3").Value = 16.55
$ws.Range("AE3").Value = 9690
$ws.Range("AF3").Value = 1.9
$ws.Range("AG3").Value = 125
$ws.Range("AH3").Value = 0.68
$ws.Range("AI3").Value = 14.58
$ws.Range("AJ3").Value = 16460000

# Row 4
$ws.Range("D4").Value = 1889
$ws.Range("E4").Value = 132
$ws.Range("F4").Value = 132
$ws.Range("G4").Value = 226
$ws.Range("H4").Value = 189
$ws.Range("I4").Value = 189
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 2400
$ws.Range("L4").Value = 642
$ws.Range("M4").Value = 1758
$ws.Range("N4").Value = 1758
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 41
$ws.Range("Q4").Value = 72
$ws.Range("R4").Value = -13
$ws.Range("S4").Value = -51
$ws.Range("T4").Value = 2
$ws.Range("U4").Value = 70
$ws.Range("V4").ClearContents()
$ws.Range("W4").Value = 6.97
$ws.Range("X4").Value = 9.98
$ws.Range("Y4").Value = 11.25
$ws.Range("Z4").Value = 8.23
$ws.Range("AA4").Value = 36.52
$ws.Range("AB4").Value = 4323.12
$ws.Range("AC4").Value = 1146
$ws.Range("AD4").Value = 13.11
$ws.Range("AE4").Value = 10682
$ws.Range("AF4").Value = 1.41
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 1
$ws.Range("AI4").Value = 13.09
$ws.Range("AJ4").Value = 16460000

# Row 5
$ws.Range("D5").Value = 2640
$ws.Range("E5").Value = 148
$ws.Range("F5").Value = 148
$ws.Range("G5").Value = 291
$ws.Range("H5").Value = 238
$ws.Range("I5").Value = 238
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 2759
$ws.Range("L5").Value = 792
$ws.Range("M5").Value = 1968
$ws.Range("N5").Value = 1968
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 82
$ws.Range("Q5").Value = 62
$ws.Range("R5").Value = 57
$ws.Range("S5").Value = -25
$ws.Range("T5").Value = 1
$ws.Range("U5").Value = 61
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 5.6
$ws.Range("X5").Value = 9
$ws.Range("Y5").Value = 12.75
$ws.Range("Z5").Value = 9.21
$ws.Range("AA5").Value = 40.23
$ws.Range("AB5").Value = 2373.51
$ws.Range("AC5").Value = 1443
$ws.Range("AD5").Value = 10.11
$ws.Range("AE5").Value = 11953
$ws.Range("AF5").Value = 1.22
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 1.71
$ws.Range("AI5").Value = 17.32
$ws.Range("AJ5").Value = 16460000

# Row 6
$ws.Range("D6").Value = 2720
$ws.Range("E6").Value = 147
$ws.Range("F6").Value = 147
$ws.Range("G6").Value = 323
$ws.Range("H6").Value = 242
$ws.Range("I6").Value = 242
$ws.Range("K6").Value = 2900
$ws.Range("L6").Value = 824
$ws.Range("M6").Value = 2076
$ws.Range("N6").Value = 2076
$ws.Range("P6").Value = 82
$ws.Range("Q6").Value = 371
$ws.Range("R6").Value = -350
$ws.Range("S6").Value = -34
$ws.Range("T6").Value = 8
$ws.Range("U6").Value = 363
$ws.Range("V6").Value = 7
$ws.Range("W6").Value = 5.4
$ws.Range("X6").Value = 8.91
$ws.Range("Y6").Value = 11.99
$ws.Range("Z6").Value = 8.57
$ws.Range("AA6").Value = 39.69
$ws.Range("AB6").Value = 2499.11
$ws.Range("AC6").Value = 1473
$ws.Range("AD6").Value = 6.46
$ws.Range("AE6").Value = 12614
$ws.Range("AF6").Value = 0.75
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 3.68
$ws.Range("AI6").Value = 23.77
$ws.Range("AJ6").Value = 16460000

# Row 7
$ws.Range("D7").Value = 2450
$ws.Range("E7").Value = 126
$ws.Range("G7").Value = 245
$ws.Range("H7").Value = 178
$ws.Range("I7").Value = 170
$ws.Range("K7").Value = 2862
$ws.Range("L7").Value = 762
$ws.Range("M7").Value = 2100
$ws.Range("N7").Value = 2100
$ws.Range("P7").Value = 81
$ws.Range("Q7").Value = 280
$ws.Range("R7").Value = -32
$ws.Range("S7").Value = -52
$ws.Range("T7").Value = 7
$ws.Range("U7").ClearContents()
$ws.Range("W7").Value = 5.14
$ws.Range("X7").Value = 7.27
$ws.Range("Y7").Value = 8.16
$ws.Range("Z7").Value = 6.18
$ws.Range("AA7").Value = 36.29
$ws.Range("AC7").Value = 1036
$ws.Range("AD7").Value = 9.08
$ws.Range("AE7").Value = 12758
$ws.Range("AF7").Value = 0.74
$ws.Range("AG7").Value = 350
$ws.Range("AH7").Value = 3.72
$ws.Range("AI7").Value = 33.79

# Row 8
$ws.Range("D8").Value = 2556
$ws.Range("E8").Value = 134
$ws.Range("G8").Value = 258
$ws.Range("H8").Value = 189
$ws.Range("I8").Value = 180
$ws.Range("K8").Value = 3036
$ws.Range("L8").Value = 816
$ws.Range("M8").Value = 2221
$ws.Range("N8").Value = 2221
$ws.Range("P8").Value = 81
$ws.Range("Q8").Value = 341
$ws.Range("R8").Value = -39
$ws.Range("S8").Value = -60
$ws.Range("T8").Value = 4
$ws.Range("U8").ClearContents()
$ws.Range("W8").Value = 5.23
$ws.Range("X8").Value = 7.39
$ws.Range("Y8").Value = 8.33
$ws.Range("Z8").Value = 6.41
$ws.Range("AA8").Value = 36.72
$ws.Range("AC8").Value = 1094
$ws.Range("AD8").Value = 8.6
$ws.Range("AE8").Value = 13493
$ws.Range("AF8").Value = 0.7
$ws.Range("AG8").Value = 350
$ws.Range("AH8").Value = 3.72
$ws.Range("AI8").Value = 32.01

# Row 9
$ws.Range("D9").Value = 2734
$ws.Range("E9").Value = 149
$ws.Range("G9").Value = 250
$ws.Range("H9").Value = 197
$ws.Range("I9").Value = 197
$ws.Range("K9").Value = 3238
$ws.Range("L9").Value = 878
$ws.Range("M9").Value = 2360
$ws.Range("N9").Value = 2360
$ws.Range("P9").Value = 81
$ws.Range("Q9").Value = 277
$ws.Range("R9").Value = -50
$ws.Range("S9").Value = -60
$ws.Range("T9").Value = 4
$ws.Range("U9").ClearContents()
$ws.Range("W9").Value = 5.45
$ws.Range("X9").Value = 7.2
$ws.Range("Y9").Value = 8.6
$ws.Range("Z9").Value = 6.28
$ws.Range("AA9").Value = 37.23
$ws.Range("AC9").Value = 1197
$ws.Range("AD9").Value = 7.86
$ws.Range("AE9").Value = 14335
$ws.Range("AF9").Value = 0.66
$ws.Range("AG9").Value = 350
$ws.Range("AH9").Value = 3.72
$ws.Range("AI9").Value = 29.24
